$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.093070387840271
$ws.Range("B1").Value = 1.005301475524902
$ws.Range("C1").Value = 15
$ws.Range("D1").Value = 1.996190309524536
$ws.Range("E1").Value = 1.014667630195618
